# Update loading_percent values for Case_5_51 (380 kV case) in res_line sheet.
# The sheet has a header row (row 1: column indices 0..13 in B1:O1) and
# data rows 2-25 (row index 0..23 in column A). Columns B-F, H-I, K-L, N
# carry the per-line loading percentages that were recomputed for this run;
# columns G, J, M, O (always 0) and column A (the index) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.14118879991443
$ws.Range("C2").Value = 4.714971814194778
$ws.Range("D2").Value = 7.602437512283195
$ws.Range("E2").Value = 9.938551248749587
$ws.Range("F2").Value = 37.75385985877477
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 31.0183022204402
$ws.Range("K2").Value = 13.98041708908117
$ws.Range("L2").Value = 10.35233306271483
$ws.Range("N2").Value = 22.18326019142633
$ws.Range("B3").Value = 15.94619127431526
$ws.Range("C3").Value = 4.457157147758807
$ws.Range("D3").Value = 7.609906557120054
$ws.Range("E3").Value = 9.947211762593216
$ws.Range("F3").Value = 37.67890564883056
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 31.04199030422727
$ws.Range("K3").Value = 13.84376505895687
$ws.Range("L3").Value = 10.34271338168329
$ws.Range("N3").Value = 22.23616157705349
$ws.Range("B4").Value = 15.82962789282467
$ws.Range("C4").Value = 4.290122461420085
$ws.Range("D4").Value = 7.614604908547721
$ws.Range("E4").Value = 9.954011870861304
$ws.Range("F4").Value = 37.64140246172956
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 31.06210622571418
$ws.Range("K4").Value = 13.76283001884176
$ws.Range("L4").Value = 10.33875354970148
$ws.Range("N4").Value = 22.27054980033538
$ws.Range("B5").Value = 15.78298076116343
$ws.Range("C5").Value = 4.219881956196169
$ws.Range("D5").Value = 7.616547982674423
$ws.Range("E5").Value = 9.957156025521193
$ws.Range("F5").Value = 37.62826945504845
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 31.07170291591033
$ws.Range("K5").Value = 13.73063072870052
$ws.Range("L5").Value = 10.33763088198543
$ws.Range("N5").Value = 22.28504317742446
$ws.Range("B6").Value = 15.77528822334052
$ws.Range("C6").Value = 4.208088178793431
$ws.Range("D6").Value = 7.616872353647087
$ws.Range("E6").Value = 9.957700647828419
$ws.Range("F6").Value = 37.62621876036369
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 31.07338089545446
$ws.Range("K6").Value = 13.7253323245896
$ws.Range("L6").Value = 10.3374741563766
$ws.Range("N6").Value = 22.28747879053432
$ws.Range("B7").Value = 15.82899526354918
$ws.Range("C7").Value = 4.289183933733484
$ws.Range("D7").Value = 7.61463099803536
$ws.Range("E7").Value = 9.954052763246946
$ws.Range("F7").Value = 37.64121663218338
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 31.06222998696243
$ws.Range("K7").Value = 13.76239255424162
$ws.Range("L7").Value = 10.3387364192484
$ws.Range("N7").Value = 22.27074331940016
$ws.Range("B8").Value = 16.07333145965082
$ws.Range("C8").Value = 4.627903070063601
$ws.Range("D8").Value = 7.604989653914703
$ws.Range("E8").Value = 9.941229876895443
$ws.Range("F8").Value = 37.72625270988098
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 31.02531265243606
$ws.Range("K8").Value = 13.93270490890872
$ws.Range("L8").Value = 10.34861313697545
$ws.Range("N8").Value = 22.20110516953982
$ws.Range("B9").Value = 16.57500271281521
$ws.Range("C9").Value = 5.222029863805464
$ws.Range("D9").Value = 7.586964278748884
$ws.Range("E9").Value = 9.927833649374593
$ws.Range("F9").Value = 37.96018208884475
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 30.99718470839053
$ws.Range("K9").Value = 14.28861576209479
$ws.Range("L9").Value = 10.38335069873549
$ws.Range("N9").Value = 22.07965062623
$ws.Range("B10").Value = 16.95377678489793
$ws.Range("C10").Value = 5.615050713781692
$ws.Range("D10").Value = 7.574243821411605
$ws.Range("E10").Value = 9.925132632703974
$ws.Range("F10").Value = 38.17231739666796
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 31.00356771817159
$ws.Range("K10").Value = 14.56119984134623
$ws.Range("L10").Value = 10.41812218228075
$ws.Range("N10").Value = 21.99959563387453
$ws.Range("B11").Value = 17.12754711217011
$ws.Range("C11").Value = 5.784314008354164
$ws.Range("D11").Value = 7.568567353477798
$ws.Range("E11").Value = 9.925448254422871
$ws.Range("F11").Value = 38.27738100427614
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 31.01234896041705
$ws.Range("K11").Value = 14.68711586200626
$ws.Range("L11").Value = 10.43591742168708
$ws.Range("N11").Value = 21.96516308156763
$ws.Range("B12").Value = 17.19349710297069
$ws.Range("C12").Value = 5.847035573848112
$ws.Range("D12").Value = 7.566433428165451
$ws.Range("E12").Value = 9.925789139298526
$ws.Range("F12").Value = 38.31837823184302
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 31.0165187554083
$ws.Range("K12").Value = 14.73503021523648
$ws.Range("L12").Value = 10.44293717139401
$ws.Range("N12").Value = 21.95240939979748
$ws.Range("B13").Value = 17.17928805913366
$ws.Range("C13").Value = 5.833588567973571
$ws.Range("D13").Value = 7.566892315299771
$ws.Range("E13").Value = 9.9257058905844
$ws.Range("F13").Value = 38.30949517181314
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 31.01558317066608
$ws.Range("K13").Value = 14.7247013346
$ws.Range("L13").Value = 10.44141289712712
$ws.Range("N13").Value = 21.95514345359132
$ws.Range("B14").Value = 17.13297027344466
$ws.Range("C14").Value = 5.789501724624587
$ws.Range("D14").Value = 7.568391482094962
$ws.Range("E14").Value = 9.925471867104722
$ws.Range("F14").Value = 38.28072969931318
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 31.0126750919212
$ws.Range("K14").Value = 14.69105338647547
$ws.Range("L14").Value = 10.436489327362
$ws.Range("N14").Value = 21.96410811440602
$ws.Range("B15").Value = 17.1046165574532
$ws.Range("C15").Value = 5.762318103348362
$ws.Range("D15").Value = 7.569311794389817
$ws.Range("E15").Value = 9.925357326765564
$ws.Range("F15").Value = 38.26326724006557
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 31.01100376411511
$ws.Range("K15").Value = 14.67047207875374
$ws.Range("L15").Value = 10.43351000086108
$ws.Range("N15").Value = 21.96963635982538
$ws.Range("B16").Value = 16.94244448914608
$ws.Range("C16").Value = 5.603796818454523
$ws.Range("D16").Value = 7.574616990365565
$ws.Range("E16").Value = 9.925143023620633
$ws.Range("F16").Value = 38.16562192713796
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 31.00311207041839
$ws.Range("K16").Value = 14.55300582779606
$ws.Range("L16").Value = 10.41699874019024
$ws.Range("N16").Value = 22.00188580146367
$ws.Range("B17").Value = 16.84328663980045
$ws.Range("C17").Value = 5.504105443785918
$ws.Range("D17").Value = 7.577899616179804
$ws.Range("E17").Value = 9.925406621661793
$ws.Range("F17").Value = 38.10789951980728
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 30.99977585215266
$ws.Range("K17").Value = 14.48140437743187
$ws.Range("L17").Value = 10.40737386360281
$ws.Range("N17").Value = 22.02217798145907
$ws.Range("B18").Value = 16.78639523560454
$ws.Range("C18").Value = 5.44586933755812
$ws.Range("D18").Value = 7.579798074119743
$ws.Range("E18").Value = 9.925703643892172
$ws.Range("F18").Value = 38.07550630683651
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 30.99841024439702
$ws.Range("K18").Value = 14.44040404724464
$ws.Range("L18").Value = 10.40202423827965
$ws.Range("N18").Value = 22.03403634459961
$ws.Range("B19").Value = 16.76715905078292
$ws.Range("C19").Value = 5.425997852851838
$ws.Range("D19").Value = 7.580442648194269
$ws.Range("E19").Value = 9.9258292044179
$ws.Range("F19").Value = 38.06467771530889
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 30.99804291422321
$ws.Range("K19").Value = 14.42655476404183
$ws.Range("L19").Value = 10.4002450479126
$ws.Range("N19").Value = 22.03808347999336
$ws.Range("B20").Value = 16.853827974039
$ws.Range("C20").Value = 5.514810565374409
$ws.Range("D20").Value = 7.577549102419227
$ws.Range("E20").Value = 9.925363516888117
$ws.Range("F20").Value = 38.11396076154885
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 31.00007373725082
$ws.Range("K20").Value = 14.48900785193391
$ws.Range("L20").Value = 10.40837918456526
$ws.Range("N20").Value = 22.01999850902004
$ws.Range("B21").Value = 17.14657144384803
$ws.Range("C21").Value = 5.802488441879781
$ws.Range("D21").Value = 7.567950717919447
$ws.Range("E21").Value = 9.925534603763483
$ws.Range("F21").Value = 38.28914608769539
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 31.01350635262741
$ws.Range("K21").Value = 14.70093063055854
$ws.Range("L21").Value = 10.43792789841698
$ws.Range("N21").Value = 21.96146723837947
$ws.Range("B22").Value = 17.33872398274552
$ws.Range("C22").Value = 5.982488270641579
$ws.Range("D22").Value = 7.56176862639069
$ws.Range("E22").Value = 9.926936314453819
$ws.Range("F22").Value = 38.41069358841626
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 31.02720724094753
$ws.Range("K22").Value = 14.84077092944439
$ws.Range("L22").Value = 10.45887630047193
$ws.Range("N22").Value = 21.92487581086609
$ws.Range("B23").Value = 17.23611355096811
$ws.Range("C23").Value = 5.887153509513069
$ws.Range("D23").Value = 7.565059865953917
$ws.Range("E23").Value = 9.926070427078294
$ws.Range("F23").Value = 38.34518282887214
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 31.01944481412901
$ws.Range("K23").Value = 14.76602736731553
$ws.Range("L23").Value = 10.44754715623637
$ws.Range("N23").Value = 21.94425333293493
$ws.Range("B24").Value = 16.84906187134824
$ws.Range("C24").Value = 5.50997365204833
$ws.Range("D24").Value = 7.577707534656808
$ws.Range("E24").Value = 9.925382551370623
$ws.Range("F24").Value = 38.11121800503228
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 30.99993734246426
$ws.Range("K24").Value = 14.48556980605158
$ws.Range("L24").Value = 10.40792410605801
$ws.Range("N24").Value = 22.02098324997944
$ws.Range("B25").Value = 16.43723882974677
$ws.Range("C25").Value = 5.068898992790525
$ws.Range("D25").Value = 7.591747772672372
$ws.Range("E25").Value = 9.930201710056439
$ws.Range("F25").Value = 37.88976931403761
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 31.00004669632435
$ws.Range("K25").Value = 14.19021587105771
$ws.Range("L25").Value = 10.37231901989337
$ws.Range("N25").Value = 22.11089343246495
